# "adding averages and more checks"
# - Training Dashboard: PERIOD TO EXPIRE (H) recomputed and LAST UPDATE (I)
#   bumped forward to 16-Sep-2025 for every training row.
# - Exam Dashboard: widen the COMMENTS column and make the remark more
#   descriptive ("OK" -> "date is valid").
# - Header / title styling: bold white text on the header band (and the
#   title now shares that same bold/white look instead of the old 14pt size).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Training Dashboard"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=3;  H=408},
    @{Row=4;  H=327},
    @{Row=5;  H=328},
    @{Row=6;  H=357},
    @{Row=7;  H=371},
    @{Row=8;  H=672},
    @{Row=9;  H=409},
    @{Row=10; H=373},
    @{Row=11; H=377},
    @{Row=12; H=407},
    @{Row=13; H=387},
    @{Row=14; H=391},
    @{Row=15; H=395},
    @{Row=16; H=133},
    @{Row=17; H=408},
    @{Row=18; H=357},
    @{Row=19; H=117},
    @{Row=20; H=174},
    @{Row=21; H=177},
    @{Row=22; H=189},
    @{Row=23; H=229}
)

foreach ($u in $updates) {
    $ws1.Cells.Item($u.Row, 8).Value = $u.H
}

# LAST UPDATE (column I) moves from 08-Sep-2025 to 16-Sep-2025 for every
# data row. Assigning a date-look-alike string straight to .Value would be
# auto-parsed into a real Excel date, so stamp it in via a literal-text
# formula and freeze it back down to a plain value with paste-special.
$lastUpdateRange = $ws1.Range("I3:I23")
$lastUpdateRange.Formula = "=""16-Sep-2025"""
$lastUpdateRange.Copy()
$lastUpdateRange.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet 2: "Exam Dashboard"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# COMMENTS column gets a bit wider.
$ws2.Columns(5).ColumnWidth = 14.17

# First remark becomes more descriptive.
$ws2.Cells.Item(3, 5).Value = "date is valid"

# ---------------------------------------------------------------------
# Header / title styling (both sheets share the same style table): the
# header band (row 2) text turns bold white, and the dashboard title
# switches to that same bold/white look.
# ---------------------------------------------------------------------
foreach ($ws in @($ws1, $ws2)) {
    $headerRow = $ws.Cells.Item(2, 1).EntireRow
    $headerRow.Font.Bold = $true
    $headerRow.Font.Color = 16777215

    $title = $ws.Cells.Item(1, 1)
    $title.Font.Bold = $true
    $title.Font.Size = 11
    $title.Font.Color = 16777215
}
